$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 1885
$ws.Range("E2").Value = 147
$ws.Range("F2").Value = 147
$ws.Range("G2").Value = 112
$ws.Range("H2").Value = 71
$ws.Range("I2").Value = 75
$ws.Range("J2").Value = -4
$ws.Range("K2").Value = 2846
$ws.Range("L2").Value = 1157
$ws.Range("M2").Value = 1689
$ws.Range("N2").Value = 1615
$ws.Range("O2").Value = 75
$ws.Range("P2").Value = 108
$ws.Range("Q2").Value = 147
$ws.Range("R2").Value = -163
$ws.Range("S2").Value = 73
$ws.Range("T2").Value = 218
$ws.Range("U2").Value = -71
$ws.Range("V2").Value = 624
$ws.Range("W2").Value = 7.78
$ws.Range("X2").Value = 3.77
$ws.Range("Y2").Value = 4.73
$ws.Range("Z2").Value = 2.63
$ws.Range("AA2").Value = 68.5
$ws.Range("AB2").Value = 1461.77
$ws.Range("AC2").Value = 347
$ws.Range("AD2").Value = 13.46
$ws.Range("AE2").Value = 9174
$ws.Range("AF2").Value = 0.51
$ws.Range("AG2").Value = 200
$ws.Range("AH2").Value = 4.28
$ws.Range("AI2").Value = 46.35
$ws.Range("AJ2").Value = 21691750

$ws.Range("D3").Value = 1964
$ws.Range("E3").Value = 121
$ws.Range("F3").Value = 121
$ws.Range("G3").Value = 87
$ws.Range("H3").Value = 70
$ws.Range("I3").Value = 68
$ws.Range("J3").Value = 2
$ws.Range("K3").Value = 3163
$ws.Range("L3").Value = 1417
$ws.Range("M3").Value = 1747
$ws.Range("N3").Value = 1663
$ws.Range("O3").Value = 83
$ws.Range("P3").Value = 108
$ws.Range("Q3").Value = -99
$ws.Range("R3").Value = -143
$ws.Range("S3").Value = 231
$ws.Range("T3").Value = 165
$ws.Range("U3").Value = -263
$ws.Range("V3").Value = 868
$ws.Range("W3").Value = 6.17
$ws.Range("X3").Value = 3.55
$ws.Range("Y3").Value = 4.13
$ws.Range("Z3").Value = 2.32
$ws.Range("AA3").Value = 81.11
$ws.Range("AB3").Value = 1492.01
$ws.Range("AC3").Value = 312
$ws.Range("AD3").Value = 16.32
$ws.Range("AE3").Value = 9310
$ws.Range("AF3").Value = 0.55
$ws.Range("AG3").Value = 240
$ws.Range("AH3").Value = 4.72
$ws.Range("AI3").Value = 63.38
$ws.Range("AJ3").Value = 21691750

$ws.Range("D4").Value = 2040
$ws.Range("E4").Value = 92
$ws.Range("F4").Value = 92
$ws.Range("G4").Value = 77
$ws.Range("H4").Value = 64
$ws.Range("I4").Value = 65
$ws.Range("J4").Value = -2
$ws.Range("K4").Value = 3209
$ws.Range("L4").Value = 1369
$ws.Range("M4").Value = 1840
$ws.Range("N4").Value = 1713
$ws.Range("O4").Value = 127
$ws.Range("P4").Value = 108
$ws.Range("Q4").Value = 160
$ws.Range("R4").Value = -141
$ws.Range("S4").Value = -74
$ws.Range("T4").Value = 154
$ws.Range("U4").Value = 6
$ws.Range("V4").Value = 769
$ws.Range("W4").Value = 4.5
$ws.Range("X4").Value = 3.12
$ws.Range("Y4").Value = 3.86
$ws.Range("Z4").Value = 2
$ws.Range("AA4").Value = 74.44
$ws.Range("AB4").Value = 1536.47
$ws.Range("AC4").Value = 301
$ws.Range("AD4").Value = 19.29
$ws.Range("AE4").Value = 9588
$ws.Range("AF4").Value = 0.6
$ws.Range("AG4").Value = 265
$ws.Range("AH4").Value = 4.57
$ws.Range("AI4").Value = 72.58
$ws.Range("AJ4").Value = 21691750

$ws.Range("D5").Value = 2065
$ws.Range("E5").Value = 84
$ws.Range("F5").Value = 84
$ws.Range("G5").Value = 129
$ws.Range("H5").Value = 101
$ws.Range("I5").Value = 102
$ws.Range("J5").Value = -1
$ws.Range("K5").Value = 3359
$ws.Range("L5").Value = 1437
$ws.Range("M5").Value = 1922
$ws.Range("N5").Value = 1769
$ws.Range("O5").Value = 152
$ws.Range("P5").Value = 108
$ws.Range("Q5").Value = 158
$ws.Range("R5").Value = -128
$ws.Range("S5").Value = 8
$ws.Range("T5").Value = 121
$ws.Range("U5").Value = 37
$ws.Range("V5").Value = 795
$ws.Range("W5").Value = 4.09
$ws.Range("X5").Value = 4.88
$ws.Range("Y5").Value = 5.83
$ws.Range("Z5").Value = 3.07
$ws.Range("AA5").Value = 74.8
$ws.Range("AB5").Value = 1589.1
$ws.Range("AC5").Value = 468
$ws.Range("AD5").Value = 12.84
$ws.Range("AE5").Value = 9905
$ws.Range("AF5").Value = 0.61
$ws.Range("AG5").Value = 295
$ws.Range("AH5").Value = 4.91
$ws.Range("AI5").Value = 51.9
$ws.Range("AJ5").Value = 21691750

$ws.Range("D6").Value = 2127
$ws.Range("E6").Value = 148
$ws.Range("F6").Value = 148
$ws.Range("G6").Value = 111
$ws.Range("H6").Value = 49
$ws.Range("I6").Value = 55
$ws.Range("K6").Value = 3573
$ws.Range("L6").Value = 1558
$ws.Range("M6").Value = 2015
$ws.Range("N6").Value = 1830
$ws.Range("P6").Value = 108
$ws.Range("Q6").Value = 106
$ws.Range("R6").Value = -128
$ws.Range("S6").Value = 16
$ws.Range("T6").Value = 157
$ws.Range("U6").Value = -51
$ws.Range("V6").Value = 846
$ws.Range("W6").Value = 6.94
$ws.Range("X6").Value = 2.29
$ws.Range("Y6").Value = 3.08
$ws.Range("Z6").Value = 1.41
$ws.Range("AA6").Value = 77.34
$ws.Range("AB6").Value = 1710.2
$ws.Range("AC6").Value = 255
$ws.Range("AD6").Value = 42.32
$ws.Range("AE6").Value = 10244
$ws.Range("AF6").Value = 1.05
$ws.Range("AG6").Value = 250
$ws.Range("AH6").Value = 2.31
$ws.Range("AI6").Value = 80.68
$ws.Range("AJ6").Value = 21691750

$ws.Range("D7:AI7").ClearContents()
$ws.Range("D8:AI8").ClearContents()
$ws.Range("D9:AI9").ClearContents()
